# Rework the "Bids" sheet to feed a new admin-dashboard view:
#  - add a small header/summary block at the top (Bid / College, and the
#    current leading bid of 89 for Oxford University)
#  - keep the first applicant's row (indrajit), but enrich it with a
#    phone number and his bid amount
#  - drop the other, now-unused applicant rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet - the old rows 2-6 are being replaced wholesale.
$ws.Cells.Clear()

# Row 1: new header row
$ws.Range("A1").Value = "Bid"
$ws.Range("B1").Value = "College"

# Row 2: headline bid figure for the dashboard
$ws.Range("A2").Value = 89
$ws.Range("B2").Value = "oxford university"

# Row 3: retained/extended applicant record (indrajit)
$ws.Range("A3").Value = "indrajit"
$ws.Range("B3").Value = "indrajit.chandra.MBA23@said.oxford.edu"

# Phone number must stay text (it has a leading context as an id, not a
# number to calculate with), so force a text format before assigning it.
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "9732253783"

$ws.Range("D3").Value = 89

# The phone-number-look-alike text in C3 is intentional; tell Excel to
# stop flagging it as "number stored as text".
$ws.Range("C3").Errors.Item(9).Ignore = $true
